$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rename the section bookmark:
#    missingcols-cond-ordergroup_over -> ordergroup_over-sum_cols_first
#
# Word bookmarks cannot be renamed by assigning .Name, so the bookmark must be
# deleted and re-added at the same Range. The original bookmark is collapsed
# at the very start of the document (position 0); this runtime mis-places a
# freshly-added bookmark's End marker into the following paragraph when the
# target Range is collapsed at absolute position 0. Work around this by
# temporarily inserting a placeholder character before position 0, adding the
# new bookmark right after it (position 1, still collapsed), and then removing
# the placeholder again so the bookmark naturally settles back at position 0.
$placeholder = $d.Range(0, 0)
$placeholder.InsertBefore("X")

$bm = $d.Bookmarks("missingcols-cond-ordergroup_over")
$bm.Delete()

$newBookmarkRange = $d.Range(1, 1)
$d.Bookmarks.Add("ordergroup_over-sum_cols_first", $newBookmarkRange)

$d.Range(0, 1).Delete()

# ---------------------------------------------------------------------------
# 2. Heading run text (6.4 ...): update the verbatim option summary.
$d.Content.Find.Execute("missing(cols cond %) order(group_over)", $true, $false, $false, $false, $false, `
  $true, 1, $false, "order(group_over) sum_cols_first", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. FirstParagraph body text rewritten:
#    "order(group_over) group columns by the over variable first, placing the
#     summary and dednominator columns together. The % option wihtin
#     missing() adds the percent of missing observations. The option per is
#     specified as well to include a percentage sign."
#    becomes
#    "order(group_over) can be combined with sum_cols_first."
#
# Done as two replacements so the pre-existing run-level formatting
# (plain vs. VerbatimChar) of the surviving text is preserved:
#   a) the long descriptive sentence -> "can be combined with"
#   b) the trailing " option wihtin ... percentage sign." tail (which spans
#      the remaining runs, including the one immediately before the lone
#      VerbatimChar "%" run) collapses down to just "."
# That leaves the original VerbatimChar "%" run in place; it is renamed to
# "sum_cols_first" in step 4 below.
$d.Content.Find.Execute("group columns by the over variable first, placing the summary and dednominator columns together. The", `
  $true, $false, $false, $false, $false, $true, 1, $false, "can be combined with", 2) | Out-Null

$d.Content.Find.Execute(" option wihtin missing() adds the percent of missing observations. The option per is specified as well to include a percentage sign.", `
  $true, $false, $false, $false, $false, $true, 1, $false, ".", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. The remaining verbatim "%" run (now immediately followed by "can be
#    combined with ... .") becomes "sum_cols_first".
$d.Content.Find.Execute("can be combined with %.", $true, $false, $false, $false, $false, `
  $true, 1, $false, "can be combined with sum_cols_first.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. SourceCode example line: refresh the placeholder column headers.
$d.Content.Find.Execute("(`"Summary 1`") (`"N 1`")  (`"Summary 2`") (`"N 2`") (`"Summary Overall`")", `
  $true, $false, $false, $false, $false, $true, 1, $false, `
  "(`"Summary 1`") (`"Missing 1`")  (`"Summary 0`") (`"Missing 0`") (`"Summary Overall`")", 2) | Out-Null
